# Add seq id to fieldsample_barcode map
# The sheet contained a stray duplicate row (fieldsample_barcode
# "MFD00457B") inserted between MFD00457 and MFD00458. Remove that row so
# the barcode sequence becomes contiguous again (MFD00457, MFD00458,
# MFD00459, ...), shifting all subsequent rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Delete()
